# Regenerate the 'K' column (column G) on Sheet1 with corrected strike-count
# values (K instead of Strike#), per the recomputed std/mean-derived s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value.
$kValues = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 2
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 3
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 2
    35 = 0
    36 = 2
    37 = 0
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 0
    45 = 0
    46 = 1
    47 = 0
    48 = 1
    49 = 3
    50 = 0
    51 = 0
    52 = 1
    55 = 1
    56 = 0
    57 = 2
    59 = 2
    60 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
